$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.439.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.75%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.426.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +7.23%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.42%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.38%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.433.48"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.31%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.536"
$ws.Range("D9").Style = "Normal"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.55"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.84%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.122"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.15%  "

# Row 12
$ws.Range("E12").Value = "  +0.42%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.013.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.13%  "

# Row 14
$ws.Range("E14").Value = "  -0.33%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000186"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +8.04%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.43%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.636.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.06%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.425.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.97%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.74%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.02%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.81%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "391.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.83%  "

# Row 23
$ws.Range("E23").Value = "  -0.07%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.539"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.10%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.81%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000107"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +21.98%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.31%  "

# Row 28
$ws.Range("E28").Value = "  +7.11%  "

# Row 29
$ws.Range("E29").Value = "  -0.06%  "

# Row 30
$ws.Range("E30").Value = "  +7.77%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.20%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +14.15%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.63%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.43%  "

# Row 35
$ws.Range("E35").Value = "  -0.07%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.75"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.57%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.50%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "158.91"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.06%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0780"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.89%  "

# Row 40
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "27.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.18%  "

# Row 41
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.47%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.935.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.71%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0321"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.76%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.767"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.80%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.62"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.35%  "

# Row 46
$ws.Range("E46").Value = "  +2.76%  "

# Row 47
$ws.Range("E47").Value = "  +10.02%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.473.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.32%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.79%  "

# Row 50
$ws.Range("B50").Value = "Bittensor"
$ws.Range("C50").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "296.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +12.68%  "

# Row 51
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.28%  "
